# 1) Rename the existing sheet "Planilha1" -> "EVC"
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "EVC"

# 2) Add the new "ESTRELAS" worksheet right after "EVC" (becomes the active tab)
$ws2 = $wb.Worksheets.Add([System.Type]::Missing, $ws1)
$ws2.Name = "ESTRELAS"

# 3) Header row values
$ws2.Range("A1").Value = "PONTUAÇÃO - ESTRELAS"
$ws2.Range("B1").Value = "PESO (0 - 1)"

# Reuse the exact same cell styles already used on the "EVC" sheet so no
# new font/style entries are introduced into styles.xml:
#   EVC!A1 -> bold 14, left aligned   (style used by header col A)
#   EVC!E1 -> bold 14                 (style used by header col B)
#   EVC!A2 -> left aligned            (style used by data col A rows 2-7)
$ws1.Range("A1").Copy()
$ws2.Range("A1").PasteSpecial(-4122)
$ws1.Range("E1").Copy()
$ws2.Range("B1").PasteSpecial(-4122)
$ws1.Range("A2").Copy()
$ws2.Range("A2:A7").PasteSpecial(-4122)

$ws2.Rows.Item(1).RowHeight = 19.05

# 4) Data rows: A2:A7 = 1..6 ; B2:B7 = (A-1)/5
for ($i = 2; $i -le 7; $i++) {
    $ws2.Cells.Item($i, 1).Value = $i - 1
}
$ws2.Range("B2").Formula = "=(A2-1)/5"
$ws2.Range("B3:B7").Formula = "=(A3-1)/5"

# 5) Column widths
$ws2.Columns.Item(1).ColumnWidth = 29.42
$ws2.Columns.Item(2).ColumnWidth = 14.6

# 6) Selection / view state matching the authored file
$ws2.Range("G10").Select()
